$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''27.430.58'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '''  +2.07%  '
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = '''1.795.65'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '''  +2.70%  '
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').Value = '''  +0.41%  '
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = '''338.81'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '''  +0.81%  '
$ws.Range('E5').Style = 'Normal'
$ws.Range('E6').Value = '''  +0.26%  '
$ws.Range('E6').Style = 'Normal'
$ws.Range('D7').Value = '''0.3801'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '''  +1.48%  '
$ws.Range('E7').Style = 'Normal'
$ws.Range('D8').Value = '''0.3458'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '''  +1.62%  '
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').Value = '''48.76'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '''  +0.57%  '
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = '''1.202'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '''  +0.93%  '
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = '''0.07526'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '''  +0.13%  '
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = '''1.002'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '''  +0.23%  '
$ws.Range('E12').Style = 'Normal'
$ws.Range('E13').Value = '''  +7.64%  '
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = '''6.481'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '''  +1.18%  '
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = '''1.793.93'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '''  +2.96%  '
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = '''7.087'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '''  +0.28%  '
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Value = '''0.00001106'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '''  +1.62%  '
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = '''0.06668'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '''  -0.84%  '
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').Value = '''84.96'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '''  +2.29%  '
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = '''1.001'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '''  +0.17%  '
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = '''6.547'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '''  +4.63%  '
$ws.Range('E21').Style = 'Normal'
$ws.Range('E22').Value = '''  +3.94%  '
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = '''27.423.08'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '''  +2.14%  '
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = '''12.54'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '''  -2.48%  '
$ws.Range('E24').Style = 'Normal'
$ws.Range('E25').Value = '''  -1.08%  '
$ws.Range('E25').Style = 'Normal'
$ws.Range('D26').Value = '''2.571'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '''  +5.95%  '
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = '''1.489'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '''  +0.50%  '
$ws.Range('E27').Style = 'Normal'
$ws.Range('D28').Value = '''21.49'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '''  +9.03%  '
$ws.Range('E28').Style = 'Normal'
$ws.Range('D29').Value = '''153.06'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '''  +0.19%  '
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').Value = '''1.998.94'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '''  +3.13%  '
$ws.Range('E30').Style = 'Normal'
$ws.Range('D31').Value = '''133.87'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '''  +0.83%  '
$ws.Range('E31').Style = 'Normal'
$ws.Range('D32').Value = '''4.050'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '''  -1.81%  '
$ws.Range('E32').Style = 'Normal'
$ws.Range('E33').Value = '''  +1.17%  '
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').Value = '''0.08719'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '''  +1.02%  '
$ws.Range('E34').Style = 'Normal'
$ws.Range('E35').Value = '''  +2.97%  '
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = '''1.656'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '''  -1.95%  '
$ws.Range('E36').Style = 'Normal'
$ws.Range('B37').Value = '''TheSandbox'
$ws.Range('B37').Style = 'Normal'
$ws.Range('C37').Value = '''https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('C37').Style = 'Normal'
$ws.Range('D37').Value = '''0.6924'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '''  +8.41%  '
$ws.Range('E37').Style = 'Normal'
$ws.Range('B38').Value = '''InternetComputer(DFINITY)'
$ws.Range('B38').Style = 'Normal'
$ws.Range('C38').Value = '''https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('C38').Style = 'Normal'
$ws.Range('D38').Value = '''5.468'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '''  +0.45%  '
$ws.Range('E38').Style = 'Normal'
$ws.Range('D39').Value = '''0.06398'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '''  +1.63%  '
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').Value = '''8.891'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '''  +3.72%  '
$ws.Range('E40').Style = 'Normal'
$ws.Range('E41').Value = '''  +1.33%  '
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').Value = '''0.02350'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '''  +0.02%  '
$ws.Range('E42').Style = 'Normal'
$ws.Range('E43').Value = '''  +4.16%  '
$ws.Range('E43').Style = 'Normal'
$ws.Range('E44').Value = '''  +1.07%  '
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = '''0.6470'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '''  +2.59%  '
$ws.Range('E45').Style = 'Normal'
$ws.Range('E46').Value = '''  +0.17%  '
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').Value = '''3.875'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '''  -1.33%  '
$ws.Range('E47').Style = 'Normal'
$ws.Range('D48').Value = '''2.140'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '''  +2.49%  '
$ws.Range('E48').Style = 'Normal'
$ws.Range('D49').Value = '''130.46'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '''  +0.55%  '
$ws.Range('E49').Style = 'Normal'
$ws.Range('D50').Value = '''0.07202'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '''  -0.41%  '
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Value = '''79.74'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '''  +1.64%  '
$ws.Range('E51').Style = 'Normal'
